$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing the existing data row down to row 4
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with header-like strings
$ws.Range("A3").Value = "Hello"
$ws.Range("B3").Value = "Bye"
$ws.Range("C3").Value = "Cya"
$ws.Range("D3").Value = "Lol"
$ws.Range("E3").Value = "Why"

# Add a new row of numeric data at row 5
$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 22
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 44
$ws.Range("E5").Value = 55

# Update the selection to match the target state
$ws.Range("E6").Select()
